# Insert a new weekly record at row 168 for "Vega Central Mapocho de Santiago - Haba",
# pushing the existing rows 168-173 down to 169-174, then populate the new row
# with the latest week's data (2021-11-16 / date serial 44516).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 168..173 down to 169..174, leaving row 168 empty for the new record.
$ws.Rows.Item(168).Insert()

# Fill in the new row 168 with this week's data.
$ws.Cells.Item(168, 1).Value = 9
$ws.Cells.Item(168, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(168, 3).Value = "Metropolitana"
$ws.Cells.Item(168, 4).Value = 44516
$ws.Cells.Item(168, 5).Value = 13
$ws.Cells.Item(168, 6).Value = 100112026
$ws.Cells.Item(168, 7).Value = "Haba"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 52
$ws.Cells.Item(168, 11).Value = 6000
$ws.Cells.Item(168, 12).Value = 7000
$ws.Cells.Item(168, 13).Value = 6500
$ws.Cells.Item(168, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(168, 15).Value = "Región Metropolitana"
$ws.Cells.Item(168, 16).Value = 260
$ws.Cells.Item(168, 17).Value = 25
$ws.Cells.Item(168, 18).Value = "Hortaliza"
